# Balance.xlsx - "Fixed share holding implementation" edit
#
# Zeroes out Player 2's Property Value, Debt Taken, and Debt value to be
# repaid (columns C, rows 5/14/15) on Sheet1, and moves the active
# selection from D11 to C16 to match where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Player 2 (column C) no longer carries these starting balances.
$ws.Range("C5").Value  = 0   # Property Value
$ws.Range("C14").Value = 0   # Debt Taken
$ws.Range("C15").Value = 0   # Debt value to be repaid

# Leave the selection where the author ended up.
$ws.Range("C16").Select()
